$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.4.0-snapshot-1"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---- Elements sheet updates ----
# Columns AK (37) and AL (38) have their header + all data swapped.
$els = $wb.Worksheets.Item("Elements")

$akHeader = $els.Range("AK1").Text
$alHeader = $els.Range("AL1").Text
$els.Range("AK1").Value = $alHeader
$els.Range("AL1").Value = $akHeader

$akRow3 = $els.Range("AK3").Text
$alRow3 = $els.Range("AL3").Text
$els.Range("AK3").Value = $alRow3
$els.Range("AL3").Value = $akRow3

$akRow5 = $els.Range("AK5").Text
$alRow5 = $els.Range("AL5").Text
$els.Range("AK5").Value = $alRow5
$els.Range("AL5").Value = $akRow5

$akRow6 = $els.Range("AK6").Text
$alRow6 = $els.Range("AL6").Text
$els.Range("AK6").Value = $alRow6
$els.Range("AL6").Value = $akRow6

# bestFit column widths follow the content swap: AK becomes the wide
# (long mapping text) column, AL becomes the narrow one.
$els.Columns.Item(37).ColumnWidth = 96.5
$els.Columns.Item(38).ColumnWidth = 24.166666666666668
